$p = $ppt.ActivePresentation

# Find the slide whose title contains the "php" typo and fix it to "JAVA"
# (commit: "Fixed Typos in Functions Lection").
# Title goes from "Параметри на функции в php ДЕмо"
#               to "Параметри на функции в JAVA ДЕмо"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.Shapes.HasTitle) {
        $title = $slide.Shapes.Title
        $tr = $title.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("php")
        if ($idx -ge 0) {
            # Replace "php" plus the single space that follows it with "JAVA "
            # (1-based Characters index), so the run structure collapses to a
            # clean "JAVA " run instead of leaving a stray separate space run.
            $start = $idx + 1
            $len = 4
            $sub = $tr.Characters($start, $len)
            $sub.Text = "JAVA "
        }
    }
}
